$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-looking-like-number Price (column D) cells to stay as text
# (matches source data which stores these as inline strings, not numbers)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.972.68"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").Value = "2.342.35"
$ws.Range("E3").Value = "  +1.35%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "306.84"
$ws.Range("E5").Value = "  -1.52%  "

$ws.Range("D6").Value = "100.95"
$ws.Range("E6").Value = "  -1.34%  "

$ws.Range("E7").Value = "  -4.45%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -3.12%  "

$ws.Range("D10").Value = "35.12"
$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("D11").Value = "52.19"
$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("E12").Value = "  -1.62%  "

$ws.Range("D13").Value = "0.113"
$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("D14").Value = "6.83"
$ws.Range("E14").Value = "  -2.30%  "

$ws.Range("D15").Value = "15.97"
$ws.Range("E15").Value = "  +6.64%  "

$ws.Range("D16").Value = "2.348.28"
$ws.Range("E16").Value = "  +1.42%  "

$ws.Range("D17").Value = "0.811"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").Value = "42.884.72"
$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("D19").Value = "6.23"
$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("E20").Value = "  -2.51%  "

$ws.Range("D21").Value = "11.74"
$ws.Range("E21").Value = "  -5.01%  "

$ws.Range("D22").Value = "68.01"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").Value = "237.06"
$ws.Range("E23").Value = "  -1.91%  "

$ws.Range("E24").Value = "  +1.11%  "

$ws.Range("E25").Value = "  -2.57%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "25.57"
$ws.Range("E27").Value = "  +3.78%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.32"
$ws.Range("E28").Value = "  +9.63%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "35.28"
$ws.Range("E29").Value = "  -4.08%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "9.38"
$ws.Range("E30").Value = "  -2.54%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "161.64"
$ws.Range("E31").Value = "  -3.49%  "

$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.15"
$ws.Range("E33").Value = "  -2.60%  "

$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "4.66"
$ws.Range("E34").Value = "  +6.92%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "2.48"
$ws.Range("E35").Value = "  -0.73%  "

$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "17.48"
$ws.Range("E36").Value = "  -0.89%  "

$ws.Range("D37").Value = "0.0729"
$ws.Range("E37").Value = "  -1.95%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.96"
$ws.Range("E38").Value = "  -3.62%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "1.88"
$ws.Range("E39").Value = "  -1.06%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.103"
$ws.Range("E40").Value = "  -3.16%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.113"
$ws.Range("E41").Value = "  -2.39%  "

$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "2.45"
$ws.Range("E42").Value = "  +5.57%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.013.75"
$ws.Range("E43").Value = "  +1.96%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0286"
$ws.Range("E44").Value = "  -0.93%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "18.82"
$ws.Range("E45").Value = "  -2.24%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "10.36"
$ws.Range("E46").Value = "  +4.47%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.94"
$ws.Range("E47").Value = "  -1.04%  "

$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").Value = "55.91"
$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "2.89"
$ws.Range("E49").Value = "  +0.93%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.565.78"
$ws.Range("E50").Value = "  +1.08%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "4.73"
$ws.Range("E51").Value = "  +3.06%  "
